$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column AA (row 1) - new date label "11-jul"
$ws.Range("AA1").Value = "11-jul"

# New column Y (previously missing in data rows) and new column AA values per row
$ws.Range("Y2").Value = 0
$ws.Range("AA2").Value = 0

$ws.Range("Y3").Value = 12.282111633415292
$ws.Range("AA3").Value = 11.899744633765254

$ws.Range("Y4").Value = 15.416875209161896
$ws.Range("AA4").Value = 19.596710872362944

$ws.Range("Y5").Value = 24.750029803912465
$ws.Range("AA5").Value = 24.968764356028654

$ws.Range("Y6").Value = 0
$ws.Range("AA6").Value = 0

$ws.Range("Y7").Value = 18.603854330999585
$ws.Range("AA7").Value = 14.813672327534999

$ws.Range("Y8").Value = 12.176750436582211
$ws.Range("AA8").Value = 9.1319086042206283

$ws.Range("Y9").Value = 21.0550955369051
$ws.Range("AA9").Value = 22.675626347262316

$ws.Range("Y10").Value = 23.66822495571671
$ws.Range("AA10").Value = 22.678933420429548

$ws.Range("Y11").Value = 11.666278549944053
$ws.Range("AA11").Value = 13.779379398022883

$ws.Range("Y12").Value = 0
$ws.Range("AA12").Value = 0

$ws.Range("Y13").Value = 13.213442385099647
$ws.Range("AA13").Value = 12.932019573330978

$ws.Range("Y14").Value = 0
$ws.Range("AA14").Value = 0

$ws.Range("Y15").Value = 0
$ws.Range("AA15").Value = 0

$ws.Range("Y16").Value = 12.887269331804175
$ws.Range("AA16").Value = 22.311574590031984

$ws.Range("Y17").Value = 0
$ws.Range("AA17").Value = 0

$ws.Range("Y18").Value = 0
$ws.Range("AA18").Value = 0

# Update the selection to match the post-edit state (AC5 selected, one column further right)
$ws.Range("AC5").Select()
